# Scheduled market-data refresh: update the computed price/profit columns
# (H, I, J, K, L, M, N) on each job sheet to the latest marketboard snapshot.
# Generated from the upstream commit's per-cell value changes.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road
$ws.Cells.Item(17, 8).Value = 696658.25  # H17 was 630500.3
$ws.Cells.Item(17, 10).Value = 778383.1  # J17 was 696658.5600000001
$ws.Cells.Item(17, 12).Value = 2335149.3  # L17 was 2089975.68
$ws.Cells.Item(17, 14).Value = -2335485.3  # N17 was -2090311.68

# Row 32: Automata for the People
$ws.Cells.Item(32, 8).Value = 5531.364  # H32 was 7292.143
$ws.Cells.Item(32, 9).Value = 3857  # I32 was 7770
$ws.Cells.Item(32, 10).Value = 6159.25  # J32 was 7212.5
$ws.Cells.Item(32, 11).Value = 3857  # K32 was 7770
$ws.Cells.Item(32, 12).Value = 6159.25  # L32 was 7212.5
$ws.Cells.Item(32, 13).Value = -3531  # M32 was -7444
$ws.Cells.Item(32, 14).Value = -6811.25  # N32 was -7864.5

# Row 40: Stuck in the Moment
$ws.Cells.Item(40, 8).Value = 10107942  # H40 was 12353535
$ws.Cells.Item(40, 9).Value = 4249.6  # I40 was 5233
$ws.Cells.Item(40, 11).Value = 4249.6  # K40 was 5233
$ws.Cells.Item(40, 13).Value = -4074.6  # M40 was -5058

# Row 112: Making Ends Meet
$ws.Cells.Item(112, 8).Value = 127523.19  # H112 was 127579.375
$ws.Cells.Item(112, 9).Value = 201599.2  # I112 was 251499
$ws.Cells.Item(112, 10).Value = 93852.27  # J112 was 86272.836
$ws.Cells.Item(112, 11).Value = 604797.6000000001  # K112 was 754497
$ws.Cells.Item(112, 12).Value = 281556.81  # L112 was 258818.508
$ws.Cells.Item(112, 13).Value = -603689.6000000001  # M112 was -753389
$ws.Cells.Item(112, 14).Value = -283772.81  # N112 was -261034.508

# Row 113: Amaro Kart
$ws.Cells.Item(113, 8).Value = 3597  # H113 was 3488.3
$ws.Cells.Item(113, 9).Value = 3329  # I113 was 3269.4285
$ws.Cells.Item(113, 11).Value = 3329  # K113 was 3269.4285
$ws.Cells.Item(113, 13).Value = -75  # M113 was -15.42849999999999

# Row 116: Growing Up
$ws.Cells.Item(116, 8).Value = 4087  # H116 was 3832.6875
$ws.Cells.Item(116, 9).Value = 4059.8333  # I116 was 3773.0715
$ws.Cells.Item(116, 11).Value = 4059.8333  # K116 was 3773.0715
$ws.Cells.Item(116, 13).Value = -617.8332999999998  # M116 was -331.0715

# Row 129: Practical Command
$ws.Cells.Item(129, 8).Value = 12173.75  # H129 was 9938.9
$ws.Cells.Item(129, 9).Value = 1889.4  # I129 was 1635.1428
$ws.Cells.Item(129, 11).Value = 5668.200000000001  # K129 was 4905.428400000001
$ws.Cells.Item(129, 13).Value = -668.2000000000007  # M129 was 94.57159999999931

# Row 138: All-night Crafting
$ws.Cells.Item(138, 8).Value = 5339.577  # H138 was 5347.222
$ws.Cells.Item(138, 10).Value = 6044.3613  # J138 was 6018.1313
$ws.Cells.Item(138, 12).Value = 18133.0839  # L138 was 18054.3939
$ws.Cells.Item(138, 14).Value = -28413.0839  # N138 was -28334.3939

# Row 140: Tome for Tradition
$ws.Cells.Item(140, 8).Value = 75275  # H140 was 65994.39999999999
$ws.Cells.Item(140, 10).Value = 75275  # J140 was 65994.39999999999
$ws.Cells.Item(140, 12).Value = 75275  # L140 was 65994.39999999999
$ws.Cells.Item(140, 14).Value = -85635  # N140 was -76354.39999999999

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Cells.Item(32, 8).Value = 1269.6666  # H32 was 1271.597
$ws.Cells.Item(32, 9).Value = 1119.7539  # I32 was 1123.9849
$ws.Cells.Item(32, 11).Value = 1119.7539  # K32 was 1123.9849
$ws.Cells.Item(32, 13).Value = -832.7538999999999  # M32 was -836.9848999999999

# Row 34: Insistent Sallets
$ws.Cells.Item(34, 8).Value = 39494.5  # H34 was 39662.668

# Row 61: Dealing with the Tough Stuff
$ws.Cells.Item(61, 8).Value = 58825692  # H61 was 50001988
$ws.Cells.Item(61, 9).Value = 71430456  # I61 was 58825256
$ws.Cells.Item(61, 11).Value = 71430456  # K61 was 58825256
$ws.Cells.Item(61, 13).Value = -71430244  # M61 was -58825044

# Row 74: As the Bolt Flies
$ws.Cells.Item(74, 8).Value = 24394818  # H74 was 25645980
$ws.Cells.Item(74, 9).Value = 28575552  # I74 was 30307576
$ws.Cells.Item(74, 11).Value = 28575552  # K74 was 30307576
$ws.Cells.Item(74, 13).Value = -28574678  # M74 was -30306702

# Row 77: Heavy Metal Banned (L)
$ws.Cells.Item(77, 8).Value = 24394818  # H77 was 25645980
$ws.Cells.Item(77, 9).Value = 28575552  # I77 was 30307576
$ws.Cells.Item(77, 11).Value = 142877760  # K77 was 151537880
$ws.Cells.Item(77, 13).Value = -142873392  # M77 was -151533512

# Row 102: Smells of Rich Tama-hagane
$ws.Cells.Item(102, 8).Value = 7144502.5  # H102 was 5557175
$ws.Cells.Item(102, 9).Value = 8335001  # I102 was 6251633
$ws.Cells.Item(102, 11).Value = 8335001  # K102 was 6251633
$ws.Cells.Item(102, 13).Value = -8333379  # M102 was -6250011

# Row 136: Metal with Mettle
$ws.Cells.Item(136, 8).Value = 58825692  # H136 was 50001988
$ws.Cells.Item(136, 9).Value = 71430456  # I136 was 58825256
$ws.Cells.Item(136, 11).Value = 214291368  # K136 was 176475768
$ws.Cells.Item(136, 13).Value = -214288818  # M136 was -176473218

$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt
$ws.Cells.Item(20, 8).Value = 1974.25  # H20 was 2009.4
$ws.Cells.Item(20, 9).Value = 1849.1428  # I20 was 1869.4
$ws.Cells.Item(20, 11).Value = 1849.1428  # K20 was 1869.4
$ws.Cells.Item(20, 13).Value = -1602.1428  # M20 was -1622.4

# Row 99: Meddle in Metal
$ws.Cells.Item(99, 8).Value = 1930.174  # H99 was 1975.9524
$ws.Cells.Item(99, 9).Value = 1822.4445  # I99 was 1869.0625
$ws.Cells.Item(99, 11).Value = 1822.4445  # K99 was 1869.0625
$ws.Cells.Item(99, 13).Value = -324.4445000000001  # M99 was -371.0625

# Row 105: Ingot to Wing It
$ws.Cells.Item(105, 8).Value = 3950  # H105 was 4000
$ws.Cells.Item(105, 9).Value = 3900  # I105 was 0
$ws.Cells.Item(105, 11).Value = 3900  # K105 was 0
$ws.Cells.Item(105, 13).Value = -2153  # M105 was None

# Row 107: The Gold Experience
$ws.Cells.Item(107, 8).Value = 61009.234  # H107 was 64731.375
$ws.Cells.Item(107, 9).Value = 2257  # I107 was 2333.875
$ws.Cells.Item(107, 10).Value = 144941  # J107 was 127128.875
$ws.Cells.Item(107, 11).Value = 2257  # K107 was 2333.875
$ws.Cells.Item(107, 12).Value = 144941  # L107 was 127128.875
$ws.Cells.Item(107, 13).Value = -337  # M107 was -413.875
$ws.Cells.Item(107, 14).Value = -148781  # N107 was -130968.875

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Cells.Item(31, 8).Value = 2935.6365  # H31 was 2854.361
$ws.Cells.Item(31, 9).Value = 2725.682  # I31 was 2581.0833
$ws.Cells.Item(31, 10).Value = 3355.5454  # J31 was 3400.9167
$ws.Cells.Item(31, 11).Value = 2725.682  # K31 was 2581.0833
$ws.Cells.Item(31, 12).Value = 3355.5454  # L31 was 3400.9167
$ws.Cells.Item(31, 13).Value = -2430.682  # M31 was -2286.0833
$ws.Cells.Item(31, 14).Value = -3945.5454  # N31 was -3990.9167

# Row 34: Armoires of the Rich and Famous
$ws.Cells.Item(34, 8).Value = 2935.6365  # H34 was 2854.361
$ws.Cells.Item(34, 9).Value = 2725.682  # I34 was 2581.0833
$ws.Cells.Item(34, 10).Value = 3355.5454  # J34 was 3400.9167
$ws.Cells.Item(34, 11).Value = 2725.682  # K34 was 2581.0833
$ws.Cells.Item(34, 12).Value = 3355.5454  # L34 was 3400.9167
$ws.Cells.Item(34, 13).Value = -2523.682  # M34 was -2379.0833
$ws.Cells.Item(34, 14).Value = -3759.5454  # N34 was -3804.9167

# Row 105: Zelkova, My Love
$ws.Cells.Item(105, 8).Value = 3335699.8  # H105 was 3335633.2
$ws.Cells.Item(105, 9).Value = 4002259.8  # I105 was 4002179.8
$ws.Cells.Item(105, 11).Value = 4002259.8  # K105 was 4002179.8
$ws.Cells.Item(105, 13).Value = -4000512.8  # M105 was -4000432.8

# Row 107: Built to Last
$ws.Cells.Item(107, 8).Value = 794177.2  # H107 was 850910.4399999999
$ws.Cells.Item(107, 9).Value = 1361134  # I107 was 1555593.1
$ws.Cells.Item(107, 10).Value = 146226.58  # J107 was 146227.72
$ws.Cells.Item(107, 11).Value = 1361134  # K107 was 1555593.1
$ws.Cells.Item(107, 12).Value = 146226.58  # L107 was 146227.72
$ws.Cells.Item(107, 13).Value = -1359214  # M107 was -1553673.1
$ws.Cells.Item(107, 14).Value = -150066.58  # N107 was -150067.72

# Row 132: Hull Lotta Damage
$ws.Cells.Item(132, 8).Value = 52632996  # H132 was 58825030
$ws.Cells.Item(132, 9).Value = 52632996  # I132 was 58825030
$ws.Cells.Item(132, 11).Value = 157898988  # K132 was 176475090
$ws.Cells.Item(132, 13).Value = -157896458  # M132 was -176472560

# Row 141: No Greater Treasure
$ws.Cells.Item(141, 8).Value = 87802.14  # H141 was 89744.39999999999
$ws.Cells.Item(141, 10).Value = 94136.2  # J141 was 94188.75
$ws.Cells.Item(141, 12).Value = 94136.2  # L141 was 94188.75
$ws.Cells.Item(141, 14).Value = -104496.2  # N141 was -104548.75

$ws = $wb.Worksheets.Item("CUL")
# Row 12: Butter Me Up
$ws.Cells.Item(12, 8).Value = 407.45456  # H12 was 396.0909
$ws.Cells.Item(12, 9).Value = 285.75  # I12 was 213
$ws.Cells.Item(12, 10).Value = 477  # J12 was 615.8
$ws.Cells.Item(12, 11).Value = 857.25  # K12 was 639
$ws.Cells.Item(12, 12).Value = 1431  # L12 was 1847.4
$ws.Cells.Item(12, 13).Value = -684.25  # M12 was -466
$ws.Cells.Item(12, 14).Value = -1777  # N12 was -2193.4

# Row 39: Bloody Good Tart, This
$ws.Cells.Item(39, 8).Value = 3184.2856  # H39 was 2548.3333
$ws.Cells.Item(39, 10).Value = 7082.5  # J39 was 7165
$ws.Cells.Item(39, 12).Value = 21247.5  # L39 was 21495
$ws.Cells.Item(39, 14).Value = -21835.5  # N39 was -22083

# Row 137: Creative Chocolate
$ws.Cells.Item(137, 8).Value = 6252640.5  # H137 was 6252641
$ws.Cells.Item(137, 9).Value = 20002244  # I137 was 16668887
$ws.Cells.Item(137, 10).Value = 2820.7273  # J137 was 2893.3
$ws.Cells.Item(137, 11).Value = 60006732  # K137 was 50006661
$ws.Cells.Item(137, 12).Value = 8462.1819  # L137 was 8679.900000000001
$ws.Cells.Item(137, 13).Value = -60001632  # M137 was -50001561
$ws.Cells.Item(137, 14).Value = -18662.1819  # N137 was -18879.9

# Row 140: Sweet, Sweet Bean Juice
$ws.Cells.Item(140, 8).Value = 2857.5  # H140 was 3766
$ws.Cells.Item(140, 9).Value = 2857.5  # I140 was 3766
$ws.Cells.Item(140, 11).Value = 8572.5  # K140 was 11298
$ws.Cells.Item(140, 13).Value = -3392.5  # M140 was -6118

$ws = $wb.Worksheets.Item("GSM")
# Row 33: Thaumaturge Is Magic
$ws.Cells.Item(33, 8).Value = 0  # H33 was 20000
$ws.Cells.Item(33, 9).Value = 0  # I33 was 20000
$ws.Cells.Item(33, 11).Value = 0  # K33 was 20000
$ws.Cells.Item(33, 13).Value = $null  # M33 was -19748

# Row 80: Needs More Prayerbell
$ws.Cells.Item(80, 8).Value = 4224.25  # H80 was 4319
$ws.Cells.Item(80, 9).Value = 4224.25  # I80 was 4319
$ws.Cells.Item(80, 11).Value = 4224.25  # K80 was 4319
$ws.Cells.Item(80, 13).Value = -3226.25  # M80 was -3321

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Cells.Item(83, 8).Value = 4224.25  # H83 was 4319
$ws.Cells.Item(83, 9).Value = 4224.25  # I83 was 4319
$ws.Cells.Item(83, 11).Value = 21121.25  # K83 was 21595
$ws.Cells.Item(83, 13).Value = -16129.25  # M83 was -16603

# Row 102: Put the Metal to the Peddle
$ws.Cells.Item(102, 8).Value = 7222.7856  # H102 was 7210.357
$ws.Cells.Item(102, 9).Value = 3932.2307  # I102 was 4094.25
$ws.Cells.Item(102, 10).Value = 50000  # J102 was 25907
$ws.Cells.Item(102, 11).Value = 3932.2307  # K102 was 4094.25
$ws.Cells.Item(102, 12).Value = 50000  # L102 was 25907
$ws.Cells.Item(102, 13).Value = -2310.2307  # M102 was -2472.25
$ws.Cells.Item(102, 14).Value = -53244  # N102 was -29151

# Row 122: Awarding Academic Excellence
$ws.Cells.Item(122, 8).Value = 8468.4  # H122 was 9318.375
$ws.Cells.Item(122, 9).Value = 5068.5  # I122 was 0
$ws.Cells.Item(122, 11).Value = 15205.5  # K122 was 0
$ws.Cells.Item(122, 13).Value = -12755.5  # M122 was None

# Row 126: Gold Rush Order
$ws.Cells.Item(126, 8).Value = 3734.3333  # H126 was 3734.8333
$ws.Cells.Item(126, 9).Value = 3714.875  # I126 was 3715.4375
$ws.Cells.Item(126, 11).Value = 11144.625  # K126 was 11146.3125
$ws.Cells.Item(126, 13).Value = -8674.625  # M126 was -8676.3125

# Row 132: On Board for Lar
$ws.Cells.Item(132, 8).Value = 4823239.5  # H132 was 5016109.5
$ws.Cells.Item(132, 9).Value = 5449574.5  # I132 was 5697214.5
$ws.Cells.Item(132, 11).Value = 16348723.5  # K132 was 17091643.5
$ws.Cells.Item(132, 13).Value = -16346193.5  # M132 was -17089113.5

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Cells.Item(7, 8).Value = 2487.25  # H7 was 2725
$ws.Cells.Item(7, 9).Value = 2487.25  # I7 was 2725
$ws.Cells.Item(7, 11).Value = 2487.25  # K7 was 2725
$ws.Cells.Item(7, 13).Value = -2375.25  # M7 was -2613

# Row 16: Saddle Sore
$ws.Cells.Item(16, 8).Value = 2046.4849  # H16 was 1989.2059
$ws.Cells.Item(16, 9).Value = 927.2  # I16 was 887.7619
$ws.Cells.Item(16, 11).Value = 927.2  # K16 was 887.7619
$ws.Cells.Item(16, 13).Value = -757.2  # M16 was -717.7619

# Row 29: Hands On
$ws.Cells.Item(29, 8).Value = 0  # H29 was 45000
$ws.Cells.Item(29, 10).Value = 0  # J29 was 45000
$ws.Cells.Item(29, 12).Value = 0  # L29 was 45000
$ws.Cells.Item(29, 14).Value = $null  # N29 was -45590

# Row 93: Hide to Go Seek
$ws.Cells.Item(93, 8).Value = 1983.3077  # H93 was 2181.8333
$ws.Cells.Item(93, 9).Value = 1026.3  # I93 was 1020.375
$ws.Cells.Item(93, 10).Value = 5173.3335  # J93 was 4504.75
$ws.Cells.Item(93, 11).Value = 1026.3  # K93 was 1020.375
$ws.Cells.Item(93, 12).Value = 5173.3335  # L93 was 4504.75
$ws.Cells.Item(93, 13).Value = 221.7  # M93 was 227.625
$ws.Cells.Item(93, 14).Value = -7669.3335  # N93 was -7000.75

# Row 100: Tiger in the Sack
$ws.Cells.Item(100, 8).Value = 13307149  # H100 was 13307150
$ws.Cells.Item(100, 9).Value = 22177738  # I100 was 22177740
$ws.Cells.Item(100, 11).Value = 22177738  # K100 was 22177740
$ws.Cells.Item(100, 13).Value = -22177197  # M100 was -22177199

# Row 126: Battered Books
$ws.Cells.Item(126, 8).Value = 2487.25  # H126 was 2725
$ws.Cells.Item(126, 9).Value = 2487.25  # I126 was 2725
$ws.Cells.Item(126, 11).Value = 7461.75  # K126 was 8175
$ws.Cells.Item(126, 13).Value = -4991.75  # M126 was -5705

# Row 132: Tenets of Tanning
$ws.Cells.Item(132, 8).Value = 25014198  # H132 was 50024000
$ws.Cells.Item(132, 9).Value = 31266500  # I132 was 62528750
$ws.Cells.Item(132, 10).Value = 4994  # J132 was 5000
$ws.Cells.Item(132, 11).Value = 93799500  # K132 was 187586250
$ws.Cells.Item(132, 12).Value = 14982  # L132 was 15000
$ws.Cells.Item(132, 13).Value = -93796970  # M132 was -187583720
$ws.Cells.Item(132, 14).Value = -20042  # N132 was -20060

$ws = $wb.Worksheets.Item("WVR")
# Row 81: Where the Dragonflies, the Net Catches
$ws.Cells.Item(81, 8).Value = 1207.4166  # H81 was 1198.909
$ws.Cells.Item(81, 9).Value = 1299  # I81 was 1232.2222
$ws.Cells.Item(81, 10).Value = 932.6667  # J81 was 1049
$ws.Cells.Item(81, 11).Value = 2598  # K81 was 2464.4444
$ws.Cells.Item(81, 12).Value = 1865.3334  # L81 was 2098
$ws.Cells.Item(81, 13).Value = -1537  # M81 was -1403.4444
$ws.Cells.Item(81, 14).Value = -3987.3334  # N81 was -4220

# Row 84: To Kill a Dragon on Nameday (L)
$ws.Cells.Item(84, 8).Value = 1207.4166  # H84 was 1198.909
$ws.Cells.Item(84, 9).Value = 1299  # I84 was 1232.2222
$ws.Cells.Item(84, 10).Value = 932.6667  # J84 was 1049
$ws.Cells.Item(84, 11).Value = 12990  # K84 was 12322.222
$ws.Cells.Item(84, 12).Value = 9326.666999999999  # L84 was 10490
$ws.Cells.Item(84, 13).Value = -7686  # M84 was -7018.222
$ws.Cells.Item(84, 14).Value = -19934.667  # N84 was -21098

# Row 132: Comfy Cabins
$ws.Cells.Item(132, 8).Value = 9094376  # H132 was 9437554
$ws.Cells.Item(132, 9).Value = 11365074  # I132 was 11629393
$ws.Cells.Item(132, 10).Value = 11584.728  # J132 was 12644.2
$ws.Cells.Item(132, 11).Value = 34095222  # K132 was 34888179
$ws.Cells.Item(132, 12).Value = 34754.18399999999  # L132 was 37932.60000000001
$ws.Cells.Item(132, 13).Value = -34092692  # M132 was -34885649
$ws.Cells.Item(132, 14).Value = -39814.18399999999  # N132 was -42992.60000000001
